# Design/theme swap.
#
# The author applied a new PowerPoint "Design" (Design tab -> Office
# Theme) to the deck. In OOXML terms this rewrites ppt/theme/theme1.xml
# (the theme used by the one-and-only slide master, i.e. the design that
# is actually painted on every slide) from the old "Integral" theme to
# the default "Office Theme" palette; the previous "Integral" theme is
# what PowerPoint then keeps around as the secondary theme part (the one
# wired to the notes master).
#
# theme1.xml (Integral) and theme2.xml (Office Theme) already share an
# identical fontScheme and fmtScheme (fills/lines/effects) -- the only
# real difference between the two theme parts is the 12 colour-scheme
# entries (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink). So the edit is
# driven entirely through the Color Scheme of the deck: every slide
# shares the one master theme, so re-pointing slide 1's theme colour
# scheme at the "Office" palette re-paints ppt/theme/theme1.xml for the
# whole deck.

function RGBVal([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$tcs.Colors(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1      000000
$tcs.Colors(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2      44546A
$tcs.Colors(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Colors(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Colors(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Colors(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Colors(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Colors(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Colors(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Colors(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink 954F72
